$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("航天发展", "航天发展", "长春高新"),
    @("北方稀土", "北方稀土", "航天发展"),
    @("风华高科", "包钢股份", "汉缆股份"),
    @("汉缆股份", "利欧股份", "包钢股份"),
    @("洲际油气", "云天化", "华胜天成"),
    @("利欧股份", "汉缆股份", "利欧股份"),
    @("包钢股份", "协鑫集成", "北方稀土"),
    @("云天化", "洲际油气", "豫能控股"),
    @("华胜天成", "风语筑", "洲际油气"),
    @("长春高新", "华胜天成", "江钨装备"),
    @("再升科技", "川金诺", "风华高科"),
    @("风语筑", "湖南白银", "风语筑"),
    @("中钨高新", "风华高科", "云天化"),
    @("川金诺", "岩山科技", "再升科技"),
    @("云南锗业", "和邦生物", "川金诺"),
    @("章源钨业", "巨力索具", "巨力索具"),
    @("岩山科技", "川发龙蟒", "岩山科技"),
    @("巨力索具", "华工科技", "掌阅科技"),
    @("协鑫集成", "中钨高新", "航天电子"),
    @("豫能控股", "章源钨业", "章源钨业")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
